$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 5045.25
$ws.Range("I7").Value = 1525
$ws.Range("K7").Value = 1525
$ws.Range("M7").Value = -1413

$ws.Range("H10").Value = 8684
$ws.Range("I10").Value = 5000
$ws.Range("K10").Value = 5000
$ws.Range("M10").Value = -4707

$ws.Range("H14").Value = 5045.25
$ws.Range("I14").Value = 1525
$ws.Range("K14").Value = 1525
$ws.Range("M14").Value = -1334

$ws.Range("H21").Value = 15000
$ws.Range("I21").Value = 15000
$ws.Range("K21").Value = 15000
$ws.Range("M21").Value = -14532

$ws.Range("H23").Value = 15000
$ws.Range("I23").Value = 15000
$ws.Range("K23").Value = 15000
$ws.Range("M23").Value = -14766

$ws.Range("H64").Value = 9719.076999999999
$ws.Range("I64").Value = 8450
$ws.Range("J64").Value = 9949.817999999999
$ws.Range("K64").Value = 8450
$ws.Range("L64").Value = 9949.817999999999
$ws.Range("M64").Value = -8202
$ws.Range("N64").Value = -10445.818

$ws.Range("H67").Value = 9719.076999999999
$ws.Range("I67").Value = 8450
$ws.Range("J67").Value = 9949.817999999999
$ws.Range("K67").Value = 8450
$ws.Range("L67").Value = 9949.817999999999
$ws.Range("M67").Value = -7592
$ws.Range("N67").Value = -11665.818

$ws.Range("H92").Value = 1050.8889
$ws.Range("I92").Value = 158.28572
$ws.Range("K92").Value = 158.28572
$ws.Range("M92").Value = 1089.71428

$ws.Range("H98").Value = 1509.9474
$ws.Range("I98").Value = 1393.4706
$ws.Range("K98").Value = 1393.4706
$ws.Range("M98").Value = 104.5293999999999

$ws.Range("H107").Value = 123.8
$ws.Range("I107").Value = 123.8
$ws.Range("K107").Value = 123.8
$ws.Range("M107").Value = 1796.2

$ws.Range("H122").Value = 1509.9474
$ws.Range("I122").Value = 1393.4706
$ws.Range("K122").Value = 4180.4118
$ws.Range("M122").Value = -1730.4118

$ws.Range("H125").Value = 2970.923
$ws.Range("I125").Value = 2692.9092
$ws.Range("K125").Value = 24236.1828
$ws.Range("M125").Value = -21776.1828

$ws.Range("H132").Value = 11843.261
$ws.Range("I132").Value = 10999.75
$ws.Range("K132").Value = 32999.25
$ws.Range("M132").Value = -30469.25

$ws.Range("H138").Value = 2075.5334
$ws.Range("I138").Value = 2044.9166
$ws.Range("K138").Value = 6134.7498
$ws.Range("M138").Value = -994.7497999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 300.5
$ws.Range("J3").Value = 406
$ws.Range("L3").Value = 406
$ws.Range("N3").Value = -636

$ws.Range("H4").Value = 155.25
$ws.Range("I4").Value = 132
$ws.Range("J4").Value = 225
$ws.Range("K4").Value = 132
$ws.Range("L4").Value = 225
$ws.Range("M4").Value = -16
$ws.Range("N4").Value = -457

$ws.Range("H88").Value = 836.7273
$ws.Range("I88").Value = 557.8570999999999
$ws.Range("J88").Value = 1324.75
$ws.Range("K88").Value = 557.8570999999999
$ws.Range("L88").Value = 1324.75
$ws.Range("M88").Value = -151.8570999999999
$ws.Range("N88").Value = -2136.75

$ws.Range("H91").Value = 836.7273
$ws.Range("I91").Value = 557.8570999999999
$ws.Range("J91").Value = 1324.75
$ws.Range("K91").Value = 557.8570999999999
$ws.Range("L91").Value = 1324.75
$ws.Range("M91").Value = 846.1429000000001
$ws.Range("N91").Value = -4132.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 753.2
$ws.Range("I29").Value = 716.5
$ws.Range("J29").Value = 900
$ws.Range("K29").Value = 716.5
$ws.Range("L29").Value = 900
$ws.Range("M29").Value = -427.5
$ws.Range("N29").Value = -1478

$ws.Range("H107").Value = 7206.5713
$ws.Range("I107").Value = 2664.6667
$ws.Range("J107").Value = 8445.272000000001
$ws.Range("K107").Value = 2664.6667
$ws.Range("L107").Value = 8445.272000000001
$ws.Range("M107").Value = -744.6667000000002
$ws.Range("N107").Value = -12285.272

$ws.Range("H134").Value = 2146.35
$ws.Range("I134").Value = 2146.35
$ws.Range("K134").Value = 6439.049999999999
$ws.Range("M134").Value = -3904.049999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 154
$ws.Range("I7").Value = 180
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 180
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = -67
$ws.Range("N7").Value = -276

$ws.Range("H16").Value = 929.7143
$ws.Range("I16").Value = 908.9231
$ws.Range("K16").Value = 908.9231
$ws.Range("M16").Value = -621.9231

$ws.Range("H22").Value = 1967.25
$ws.Range("I22").Value = 1967.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1967.25
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1617.25
$ws.Range("N22").ClearContents()

$ws.Range("H69").Value = 12000
$ws.Range("I69").Value = 7666.6665
$ws.Range("J69").Value = 25000
$ws.Range("K69").Value = 7666.6665
$ws.Range("L69").Value = 25000
$ws.Range("M69").Value = -6917.6665
$ws.Range("N69").Value = -26498

$ws.Range("H72").Value = 12000
$ws.Range("I72").Value = 7666.6665
$ws.Range("J72").Value = 25000
$ws.Range("K72").Value = 22999.9995
$ws.Range("L72").Value = 75000
$ws.Range("M72").Value = -19255.9995
$ws.Range("N72").Value = -82488

$ws.Range("H113").Value = 929.7143
$ws.Range("I113").Value = 908.9231
$ws.Range("K113").Value = 908.9231
$ws.Range("M113").Value = 1261.0769

$ws.Range("H132").Value = 4580.853
$ws.Range("I132").Value = 3772.3044
$ws.Range("K132").Value = 11316.9132
$ws.Range("M132").Value = -8786.913199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 38.944443
$ws.Range("I2").Value = 40.923077
$ws.Range("J2").Value = 33.8
$ws.Range("K2").Value = 245.538462
$ws.Range("L2").Value = 202.8
$ws.Range("M2").Value = -132.538462
$ws.Range("N2").Value = -428.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 80.25
$ws.Range("I2").Value = 76.3
$ws.Range("K2").Value = 76.3
$ws.Range("M2").Value = 36.7

$ws.Range("H80").Value = 2412.25
$ws.Range("I80").Value = 2004
$ws.Range("K80").Value = 2004
$ws.Range("M80").Value = -1006

$ws.Range("H83").Value = 2412.25
$ws.Range("I83").Value = 2004
$ws.Range("K83").Value = 10020
$ws.Range("M83").Value = -5028

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8582.25
$ws.Range("I7").Value = 7734.4
$ws.Range("K7").Value = 7734.4
$ws.Range("M7").Value = -7622.4

$ws.Range("H40").Value = 4985.5
$ws.Range("I40").Value = 4004.25
$ws.Range("K40").Value = 4004.25
$ws.Range("M40").Value = -3868.25

$ws.Range("H61").Value = 7978
$ws.Range("I61").Value = 7940
$ws.Range("K61").Value = 7940
$ws.Range("M61").Value = -7738

$ws.Range("H68").Value = 9754.4
$ws.Range("I68").Value = 9431.666999999999
$ws.Range("K68").Value = 9431.666999999999
$ws.Range("M68").Value = -8682.666999999999

$ws.Range("H71").Value = 9754.4
$ws.Range("I71").Value = 9431.666999999999
$ws.Range("K71").Value = 47158.335
$ws.Range("M71").Value = -43414.335

$ws.Range("H100").Value = 5767.5
$ws.Range("I100").Value = 2312.625
$ws.Range("K100").Value = 2312.625
$ws.Range("M100").Value = -1771.625

$ws.Range("H113").Value = 7978
$ws.Range("I113").Value = 7940
$ws.Range("K113").Value = 7940
$ws.Range("M113").Value = -5770

$ws.Range("H122").Value = 7099.6
$ws.Range("J122").Value = 7999.5
$ws.Range("L122").Value = 23998.5
$ws.Range("N122").Value = -28898.5

$ws.Range("H126").Value = 8582.25
$ws.Range("I126").Value = 7734.4
$ws.Range("K126").Value = 23203.2
$ws.Range("M126").Value = -20733.2

$ws.Range("H132").Value = 3599.077
$ws.Range("I132").Value = 3617.0908
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 10851.2724
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -8321.2724
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 973
$ws.Range("I100").Value = 1222.6666
$ws.Range("J100").Value = 598.5
$ws.Range("K100").Value = 2445.3332
$ws.Range("L100").Value = 1197
$ws.Range("M100").Value = -1904.3332
$ws.Range("N100").Value = -2279

$ws.Range("H132").Value = 3216.6428
$ws.Range("I132").Value = 2276.111
$ws.Range("J132").Value = 4909.6
$ws.Range("K132").Value = 6828.333
$ws.Range("L132").Value = 14728.8
$ws.Range("M132").Value = -4298.333
$ws.Range("N132").Value = -19788.8
